$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New shared strings must be interned in this exact order (matches target
# sharedStrings.xml append order): CPFL-PIRATINING, Neoenergia PE,
# Neoenergia Brasília, CPFL JAGUARI. Write the column-A cells that
# introduce brand-new text first, in that sequence, before touching
# anything else.
$ws.Cells.Item(125, 1).Value = "CPFL-PIRATINING"
$ws.Cells.Item(126, 1).Value = "Neoenergia PE"
$ws.Cells.Item(127, 1).Value = "Neoenergia Brasília"
$ws.Cells.Item(122, 1).Value = "CPFL JAGUARI"

# Remaining cells reuse already-interned strings, so ordering among them
# doesn't affect the shared-string table.
$ws.Cells.Item(122, 2).Value = "CPFL SANTA CRUZ"
$ws.Cells.Item(123, 1).Value = "EQUATORIAL PI"
$ws.Cells.Item(123, 2).Value = "EQUATORIAL PI"
$ws.Cells.Item(124, 1).Value = "ERO"
$ws.Cells.Item(124, 2).Value = "ERO"
$ws.Cells.Item(125, 2).Value = "CPFL PIRATININGA"
$ws.Cells.Item(126, 2).Value = "CELPE"
$ws.Cells.Item(127, 2).Value = "CEB"

$ws.Range("B126").Select()

# Scroll the window so row 109 is the top visible row (mirrors
# sheetView/@topLeftCell="A109" in the target file), and put the cursor
# on B126 to match the saved selection/@activeCell.
$win = $excel.ActiveWindow
$win.ScrollRow = 109
$win.ScrollColumn = 1
